# Auto-generated: update Leve profit-calculation sheets with refreshed market-price data
# (per scheduled market-data runner) across ALC/ARM/BSM/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether
$ws.Range("H15").Value = 94.86
$ws.Range("I15").Value = 94.86
$ws.Range("K15").Value = 284.58
$ws.Range("M15").Value = -115.58

# Row 43: Growing Is Knowing
$ws.Range("H43").Value = 2376.4
$ws.Range("I43").Value = 2980
$ws.Range("J43").Value = 2225.5
$ws.Range("K43").Value = 2980
$ws.Range("L43").Value = 2225.5
$ws.Range("M43").Value = -2911
$ws.Range("N43").Value = -2363.5

# Row 74: Adhesive of Antipathy
$ws.Range("H74").Value = 50005500
$ws.Range("I74").Value = 100000000
$ws.Range("J74").Value = 11000
$ws.Range("K74").Value = 100000000
$ws.Range("L74").Value = 11000
$ws.Range("M74").Value = -99999064
$ws.Range("N74").Value = -12872

# Row 76: Warding Off Temptation
$ws.Range("H76").Value = 3725
$ws.Range("I76").Value = 3300
$ws.Range("K76").Value = 3300
$ws.Range("M76").Value = -2985

# Row 77: It's Gonna Grow Back (L)
$ws.Range("H77").Value = 50005500
$ws.Range("I77").Value = 100000000
$ws.Range("J77").Value = 11000
$ws.Range("K77").Value = 500000000
$ws.Range("L77").Value = 55000
$ws.Range("M77").Value = -499995320
$ws.Range("N77").Value = -64360

# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 3725
$ws.Range("I79").Value = 3300
$ws.Range("K79").Value = 3300
$ws.Range("M79").Value = -2208

# Row 141: Remedy for Reason
$ws.Range("H141").Value = 2911.3333
$ws.Range("I141").Value = 2402.375
$ws.Range("J141").Value = 4540
$ws.Range("K141").Value = 7207.125
$ws.Range("L141").Value = 13620
$ws.Range("M141").Value = -2027.125
$ws.Range("N141").Value = -23980


$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 927.4545000000001
$ws.Range("I2").Value = 990.2857
$ws.Range("J2").Value = 817.5
$ws.Range("K2").Value = 990.2857
$ws.Range("L2").Value = 817.5
$ws.Range("M2").Value = -877.2857
$ws.Range("N2").Value = -1043.5

# Row 32: Ingot We Trust
$ws.Range("H32").Value = 7111.0967
$ws.Range("I32").Value = 5057.8486
$ws.Range("J32").Value = 12130.148
$ws.Range("K32").Value = 5057.8486
$ws.Range("L32").Value = 12130.148
$ws.Range("M32").Value = -4770.8486
$ws.Range("N32").Value = -12704.148

# Row 63: Rivets Run through It
$ws.Range("H63").Value = 8660844
$ws.Range("I63").Value = 23087486
$ws.Range("K63").Value = 23087486
$ws.Range("M63").Value = -23086800

# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 8660844
$ws.Range("I66").Value = 23087486
$ws.Range("K66").Value = 115437430
$ws.Range("M66").Value = -115433998

# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 1659
$ws.Range("I110").Value = 1789.0588
$ws.Range("K110").Value = 1789.0588
$ws.Range("M110").Value = 255.9412

# Row 116: No Scope
$ws.Range("H116").Value = 927.4545000000001
$ws.Range("I116").Value = 990.2857
$ws.Range("J116").Value = 817.5
$ws.Range("K116").Value = 990.2857
$ws.Range("L116").Value = 817.5
$ws.Range("M116").Value = 1303.7143
$ws.Range("N116").Value = -5405.5

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 2014.7922
$ws.Range("I132").Value = 1458.2273
$ws.Range("J132").Value = 5354.1816
$ws.Range("K132").Value = 4374.6819
$ws.Range("L132").Value = 16062.5448
$ws.Range("M132").Value = -1844.6819
$ws.Range("N132").Value = -21122.5448

# Row 137: Odd Instruments
$ws.Range("H137").Value = 53780
$ws.Range("J137").Value = 53780
$ws.Range("L137").Value = 53780
$ws.Range("N137").Value = -63980


$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 927.4545000000001
$ws.Range("I3").Value = 990.2857
$ws.Range("J3").Value = 817.5
$ws.Range("K3").Value = 990.2857
$ws.Range("L3").Value = 817.5
$ws.Range("M3").Value = -876.2857
$ws.Range("N3").Value = -1045.5

# Row 133: Paring Is Caring
$ws.Range("H133").Value = 50500
$ws.Range("J133").Value = 50500
$ws.Range("L133").Value = 50500
$ws.Range("N133").Value = -60620

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 1920.4832
$ws.Range("I134").Value = 1125.2236
$ws.Range("J134").Value = 6569.6924
$ws.Range("K134").Value = 3375.6708
$ws.Range("L134").Value = 19709.0772
$ws.Range("M134").Value = -840.6707999999999
$ws.Range("N134").Value = -24779.0772

# Row 137: Dagger Swagger
$ws.Range("H137").Value = 32951.11
$ws.Range("J137").Value = 32951.11
$ws.Range("L137").Value = 32951.11
$ws.Range("N137").Value = -43151.11


$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap
$ws.Range("H5").Value = 1194.975
$ws.Range("I5").Value = 351.4
$ws.Range("J5").Value = 7100
$ws.Range("K5").Value = 1054.2
$ws.Range("L5").Value = 21300
$ws.Range("M5").Value = -942.1999999999998
$ws.Range("N5").Value = -21524

# Row 95: Soup for the Soul
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 585.6607
$ws.Range("I113").Value = 589.2368
$ws.Range("J113").Value = 578.1111
$ws.Range("K113").Value = 1767.7104
$ws.Range("L113").Value = 1734.3333
$ws.Range("M113").Value = 402.2896000000001
$ws.Range("N113").Value = -6074.3333

# Row 117: A Good Omen
$ws.Range("H117").Value = 7225.7
$ws.Range("I117").Value = 5057
$ws.Range("J117").Value = 7466.6665
$ws.Range("K117").Value = 15171
$ws.Range("L117").Value = 22399.9995
$ws.Range("M117").Value = -11729
$ws.Range("N117").Value = -29283.9995

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 8475585
$ws.Range("J131").Value = 892.2406999999999
$ws.Range("L131").Value = 2676.7221
$ws.Range("N131").Value = -12756.7221

# Row 132: More Mezcal
$ws.Range("H132").Value = 3182.8965
$ws.Range("J132").Value = 7043.4546
$ws.Range("L132").Value = 63391.0914
$ws.Range("N132").Value = -68451.0914

# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 1194.975
$ws.Range("I135").Value = 351.4
$ws.Range("J135").Value = 7100
$ws.Range("K135").Value = 3162.6
$ws.Range("L135").Value = 63900
$ws.Range("M135").Value = -627.5999999999999
$ws.Range("N135").Value = -68970

# Row 136: Simple Is Hardest
$ws.Range("H136").Value = 4380
$ws.Range("I136").Value = 4440
$ws.Range("J136").Value = 4200
$ws.Range("K136").Value = 13320
$ws.Range("L136").Value = 12600
$ws.Range("M136").Value = -8220
$ws.Range("N136").Value = -22800

# Row 137: Creative Chocolate
$ws.Range("H137").Value = 5388.4443
$ws.Range("I137").Value = 4007.5
$ws.Range("J137").Value = 6493.2
$ws.Range("K137").Value = 12022.5
$ws.Range("L137").Value = 19479.6
$ws.Range("M137").Value = -6922.5
$ws.Range("N137").Value = -29679.6


$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 6629.8667
$ws.Range("I70").Value = 6264.55
$ws.Range("J70").Value = 7360.5
$ws.Range("K70").Value = 6264.55
$ws.Range("L70").Value = 7360.5
$ws.Range("M70").Value = -5994.55
$ws.Range("N70").Value = -7900.5

# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 6629.8667
$ws.Range("I73").Value = 6264.55
$ws.Range("J73").Value = 7360.5
$ws.Range("K73").Value = 6264.55
$ws.Range("L73").Value = 7360.5
$ws.Range("M73").Value = -5328.55
$ws.Range("N73").Value = -9232.5


$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad
$ws.Range("H40").Value = 6120.4
$ws.Range("I40").Value = 6118.727
$ws.Range("K40").Value = 6118.727
$ws.Range("M40").Value = -5982.727

# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 2420.3333
$ws.Range("I93").Value = 1945.75
$ws.Range("J93").Value = 2800
$ws.Range("K93").Value = 1945.75
$ws.Range("L93").Value = 2800
$ws.Range("M93").Value = -697.75
$ws.Range("N93").Value = -5296

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 5147.7754
$ws.Range("I132").Value = 1703.65
$ws.Range("J132").Value = 7523.0347
$ws.Range("K132").Value = 5110.950000000001
$ws.Range("L132").Value = 22569.1041
$ws.Range("M132").Value = -2580.950000000001
$ws.Range("N132").Value = -27629.1041


$ws = $wb.Worksheets.Item("WVR")
# Row 49: A Leg Up on the Cold
$ws.Range("H49").Value = 33358666
$ws.Range("I49").Value = 100000000
$ws.Range("J49").Value = 38000
$ws.Range("K49").Value = 100000000
$ws.Range("L49").Value = 38000
$ws.Range("M49").Value = -99999770
$ws.Range("N49").Value = -38460

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 5651252
$ws.Range("I132").Value = 553.92
$ws.Range("J132").Value = 9806177
$ws.Range("K132").Value = 1661.76
$ws.Range("L132").Value = 29418531
$ws.Range("M132").Value = 868.2400000000002
$ws.Range("N132").Value = -29423591

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 2728.4883
$ws.Range("I136").Value = 824.8889
$ws.Range("K136").Value = 2474.6667
$ws.Range("M136").Value = 75.33329999999978

